$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13 - Pedro / 2064 / Vanessa do Romualdo / pending battery swap note
$ws.Cells.Item(13, 1).Value = "Pedro"
$ws.Cells.Item(13, 2).Value = "'2064"
$ws.Cells.Item(13, 3).Value = "Vanessa do Romualdo"
$ws.Cells.Item(13, 4).Value = "Ficou pendente a troca de bateria de alguns sensores."
$ws.Cells.Item(13, 7).Value = "Pendente"

# Row 14 - Pedro / 2261 / Escola CNA / restore comms note
$ws.Cells.Item(14, 1).Value = "Pedro"
$ws.Cells.Item(14, 2).Value = "'2261"
$ws.Cells.Item(14, 3).Value = "Escola CNA"
$ws.Cells.Item(14, 4).Value = "Tentar resutarar a comunicação da central conosco."
$ws.Cells.Item(14, 7).Value = "Pendente"

# Row 15 - Pedro / 2693 / Casa da Kenia / no beep after arming note
$ws.Cells.Item(15, 1).Value = "Pedro"
$ws.Cells.Item(15, 2).Value = "'2693"
$ws.Cells.Item(15, 3).Value = "Casa da Kenia"
$ws.Cells.Item(15, 4).Value = "Cliente disse que não está bipando após o arme."
$ws.Cells.Item(15, 7).Value = "Pendente"

# Match the author's final on-screen selection/scroll position
$ws.Range("G15").Select()
